$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny floating point precision difference on A13 (re-write value)
$ws.Cells.Item(13, 1).Value = 45863.83357263889

# Add new row 14 with sensor data
$ws.Cells.Item(14, 1).Value = 45863.87523748419
$ws.Cells.Item(14, 2).Value = 2025
$ws.Cells.Item(14, 3).Value = 30
$ws.Cells.Item(14, 4).Value = 13.33
$ws.Cells.Item(14, 5).Value = 88.92
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 6.98
$ws.Cells.Item(14, 8).Value = "ESE"
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = "21:00:20"

# Apply the same date/time number format used by the other "Fecha" cells (style s="2")
$ws.Cells.Item(14, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
